# Update crypto price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.209.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.08%  '

$ws.Range("D3").Value = "'3.417.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.05%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").Value = "'563.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.54%  '

$ws.Range("D6").Value = "'177.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.56%  '

$ws.Range("E7").Value = '  +3.46%  '

$ws.Range("D8").Value = "'3.408.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.92%  '

$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = "'0.638"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.89%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = "'0.166"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.63%  '

$ws.Range("D12").Value = "'54.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.71%  '

$ws.Range("D13").Value = "'0.0000278"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.50%  '

$ws.Range("D14").Value = "'9.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.55%  '

$ws.Range("D15").Value = "'3.956.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.27%  '

$ws.Range("D16").Value = "'18.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.80%  '

$ws.Range("D17").Value = "'3.406.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.13%  '

$ws.Range("D18").Value = "'0.118"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.46%  '

$ws.Range("E19").Value = '  +3.95%  '

$ws.Range("D20").Value = "'65.172.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.33%  '

$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("D22").Value = "'471.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.84%  '

$ws.Range("E23").Value = '  +15.31%  '

$ws.Range("D24").Value = "'4.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.72%  '

$ws.Range("D25").Value = "'87.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.11%  '

$ws.Range("D26").Value = "'13.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.43%  '

$ws.Range("D27").Value = "'10.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.45%  '

$ws.Range("D28").Value = "'2.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.29%  '

$ws.Range("D29").Value = "'8.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.30%  '

$ws.Range("D30").Value = "'30.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.33%  '

$ws.Range("D31").Value = "'6.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.94%  '

$ws.Range("D32").Value = "'11.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.29%  '

$ws.Range("D33").Value = "'581.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("E34").Value = '  +4.36%  '

$ws.Range("D35").Value = "'60.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.89%  '

$ws.Range("E36").Value = '  -0.33%  '

$ws.Range("E37").Value = '  -4.54%  '

$ws.Range("D38").Value = "'36.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.66%  '

$ws.Range("D39").Value = "'0.0₃0765"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.34%  '

$ws.Range("D40").Value = "'3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.35%  '

$ws.Range("E41").Value = '  +2.81%  '

$ws.Range("D42").Value = "'3.121.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.44%  '

$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("D44").Value = "'2.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.27%  '

$ws.Range("E45").Value = '  +3.15%  '

$ws.Range("D46").Value = "'0.0415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.92%  '

$ws.Range("D47").Value = "'3.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.97%  '

$ws.Range("E48").Value = '  +5.34%  '

$ws.Range("D49").Value = "'2.58"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'8.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.26%  '

$ws.Range("D51").Value = "'136.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.23%  '
